$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("5002825-34.2018.8.21.0002", "9000482-94.2018.8.21.0002", "Migrado"),
    @("5034684-18.2011.8.21.0001", "0111143-49.2011.8.21.0001", "Digitalizado"),
    @("5029337-62.2015.8.21.0001", "0090773-10.2015.8.21.0001", "Digitalizado"),
    @("5039775-79.2017.8.21.0001", "0039384-15.2017.8.21.0001", "Digitalizado"),
    @("5029499-57.2015.8.21.0001", "0136341-49.2015.8.21.0001", "Digitalizado"),
    @("5028728-21.2011.8.21.0001", "0344069-02.2011.8.21.0001", "Digitalizado"),
    @("5028728-21.2011.8.21.0001", "0344069-02.2011.8.21.0001", "Digitalizado"),
    @("5001806-14.2020.8.21.0134", "9000364-42.2020.8.21.0134", "Migrado"),
    @("5032111-07.2011.8.21.0001", "0111095-90.2011.8.21.0001", "Digitalizado"),
    @("5033921-46.2013.8.21.0001", "0037413-34.2013.8.21.0001", "Digitalizado"),
    @("5033938-82.2013.8.21.0001", "0084903-52.2013.8.21.0001", "Digitalizado"),
    @("5033938-82.2013.8.21.0001", "0084903-52.2013.8.21.0001", "Digitalizado"),
    @("5034971-78.2011.8.21.0001", "0108702-95.2011.8.21.0001", "Digitalizado"),
    @("5034971-78.2011.8.21.0001", "0108702-95.2011.8.21.0001", "Digitalizado"),
    @("5039859-80.2017.8.21.0001", "0129049-42.2017.8.21.0001", "Digitalizado"),
    @("5015478-52.2010.8.21.0001", "2626731-39.2010.8.21.0001", "Digitalizado"),
    @("5015600-65.2010.8.21.0001", "0155171-39.2010.8.21.0001", "Digitalizado"),
    @("5034891-17.2011.8.21.0001", "0243442-87.2011.8.21.0001", "Digitalizado"),
    @("5034891-17.2011.8.21.0001", "0243442-87.2011.8.21.0001", "Digitalizado"),
    @("5001001-39.2018.8.21.0067", "9000812-90.2018.8.21.0067", "Migrado"),
    @("5002047-92.2020.8.21.0067", "9000436-36.2020.8.21.0067", "Migrado"),
    @("5002099-25.2019.8.21.0067", "9000760-60.2019.8.21.0067", "Migrado"),
    @("5001002-24.2018.8.21.0067", "9000850-05.2018.8.21.0067", "Migrado"),
    @("5002100-10.2019.8.21.0067", "9000786-58.2019.8.21.0067", "Migrado"),
    @("5002101-92.2019.8.21.0067", "9000444-47.2019.8.21.0067", "Migrado"),
    @("5002102-77.2019.8.21.0067", "9000756-23.2019.8.21.0067", "Migrado"),
    @("5002103-62.2019.8.21.0067", "9000754-53.2019.8.21.0067", "Migrado"),
    @("5002104-47.2019.8.21.0067", "9000700-87.2019.8.21.0067", "Migrado"),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

